# Edit: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Replace the employee arrears (mora) table with an updated/expanded dataset
# and update the summary totals (VALOR MORA, Cant. Trabajadores, Cant. Periodos).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room: the table currently holds 3 data rows (16:18) followed by
#    a gap (19:22) and the signature footer (23:24). The new data needs
#    17 rows (16:32), keeping the same 4-row gap before the footer
#    (which ends up at 37:38). Insert 14 blank rows right after row 18.
# ---------------------------------------------------------------------
$ws.Rows("19:32").Insert()

# Preserve the distinctive "bottom border" look of the last table row:
# copy it from the old last row (18, still carrying that formatting)
# down onto the new last row (32) before row 18 becomes a normal row.
$ws.Range("B18:J18").Copy()
$ws.Range("B32:J32").PasteSpecial(-4122)  # xlPasteFormats

# Now give every other data row (16:31) the regular "middle" row format,
# copied from row 17 (a normal, non-bordered row).
$ws.Range("B17:J17").Copy()
$ws.Range("B16:J31").PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Fill in the new employee arrears data (rows 16-32).
# ---------------------------------------------------------------------
$rows = @(
    @("CC","73198932","ALEXANDER GUZMAN RESTREPO","2507",29120,689455),
    @("CC","73198932","ALEXANDER GUZMAN RESTREPO","2506",29120,689455),
    @("CC","73198932","ALEXANDER GUZMAN RESTREPO","2505",29120,689455),
    @("CC","73198932","ALEXANDER GUZMAN RESTREPO","2504",29120,689455),
    @("CC","73198932","ALEXANDER GUZMAN RESTREPO","2503",29120,689455),
    @("CC","73198932","ALEXANDER GUZMAN RESTREPO","2502",29120,689455),
    @("CC","1047456680","FORNARIS TORRES LUNA","2402",41600,1300000),
    @("CC","2759045","GUSTAVO MIGUEL VILLADIEGO RAMOS","2507",1898,1423500),
    @("CC","73576500","MANUEL CABARCAS ACEVEDO","2006",35112,877803),
    @("CC","1001979504","GABRIEL ISAAC LORDUY DIAZ","2506",11388,1423500),
    @("CC","1044932336","LEONAR IRIARTE DIAZ","2507",56940,1423500),
    @("CC","1050276279","CAMILA ANDREA SUAREZ TRUJILLO","2507",64000,1600000),
    @("CC","1051888670","MAURO DAVID ALMANZA PAJARO","2506",11388,1423500),
    @("CC","1143397397","MOISES TAPIAS MEZA","2506",11388,1423500),
    @("CC","20367871","ISAIAS MARIMON MARQUEZ","2205",22666,1000000),
    @("CC","1143390745","JESUS DAVID PEÑA CABRERA","2506",11388,1423500),
    @("CC","1044909687","MATEO DAVID ORTIZ RINCON","2506",11388,1423500)
)

$r = 16
foreach ($rec in $rows) {
    $ws.Cells.Item($r, 2).Value = $rec[0]   # B: Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $rec[1]   # C: N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $rec[2]   # D: Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $rec[3]   # E: Periodo Mora
    $ws.Cells.Item($r, 6).Value = $rec[4]   # F: Valor Mora
    $ws.Cells.Item($r, 7).Value = $rec[5]   # G: Salario Basico
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3) Update the summary block above the table.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 453876   # VALOR MORA
$ws.Range("C13").Value = 12       # Cant. Trabajadores
$ws.Range("F13").Value = 9        # Cant. Periodos

# ---------------------------------------------------------------------
# 4) Resize column D (Nombre Trabajador) to fit the longest new name.
# ---------------------------------------------------------------------
$ws.Columns("D").AutoFit()
